$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Mark reviewed rows with "x" in a new column N across sheets 2-5.
#    (sheet indices are 1-based; names kept here as comments for clarity)
# ---------------------------------------------------------------------------

# Sheet 2 = "sdef-ndef vs sdef-nfse"
$ws2 = $wb.Worksheets.Item(2)
foreach ($r in 1..9) {
    $ws2.Cells.Item($r, 14).Value = "x"
}
foreach ($r in 11..13) {
    $ws2.Cells.Item($r, 14).Value = "x"
}

# Sheet 3 = "sdef-ndef vs sfse-nfse"
$ws3 = $wb.Worksheets.Item(3)
foreach ($r in 11..15) {
    $ws3.Cells.Item($r, 14).Value = "x"
}

# Sheet 4 = "sdef-ndef vs sfse-ndef"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(11, 14).Value = "x"

# Sheet 5 = "sdef-nfse vs sfse-nfse"
$ws5 = $wb.Worksheets.Item(5)
foreach ($r in 1..3) {
    $ws5.Cells.Item($r, 14).Value = "x"
}

# ---------------------------------------------------------------------------
# 2) Narrow the two "long text" columns (E and K) on sheets 2-5 now that a
#    short flag column exists; this matches the new (smaller, non-bestFit)
#    column widths from the edit. The engine stores widths on a 1/6-character
#    pixel grid, so we pick the ColumnWidth input that lands closest to the
#    desired stored width.
# ---------------------------------------------------------------------------

$ws2.Columns.Item(5).ColumnWidth = 29.5
$ws2.Columns.Item(11).ColumnWidth = 31

$ws3.Columns.Item(5).ColumnWidth = 35.833333333333336
$ws3.Columns.Item(6).ColumnWidth = 3.3333333333333335
$ws3.Columns.Item(11).ColumnWidth = 42.5
$ws3.Columns.Item(12).ColumnWidth = 3.3333333333333335

$ws4.Columns.Item(5).ColumnWidth = 32.833333333333336
$ws4.Columns.Item(11).ColumnWidth = 40.5

$ws5.Columns.Item(5).ColumnWidth = 40.333333333333336
$ws5.Columns.Item(11).ColumnWidth = 44.666666666666664
$ws5.Columns.Item(12).ColumnWidth = 3.3333333333333335
$ws5.Columns.Item(13).ColumnWidth = 3.3333333333333335

# ---------------------------------------------------------------------------
# 3) Restore/update each sheet's selection (this also drives which sheet
#    ends up as the active tab: whichever sheet we touch last). Sheet 4
#    ("sdef-ndef vs sfse-ndef") must be the final active tab, so its
#    selection is applied last.
# ---------------------------------------------------------------------------

$ws2.Range("N13").Select()
$ws3.Range("N15").Select()
$ws5.Range("N6").Select()
$ws4.Range("H25").Select()

$wb.Save()
